$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "me.mods"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "com.mods"
$ws2.Columns.Item(1).ColumnWidth = 15.42578125
$ws2.Columns.Item(2).ColumnWidth = 8.85546875
Write-Host "col1 width: $($ws2.Columns.Item(1).ColumnWidth)"
